$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 9,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 12.07002266666666
$data[0,3] = 36.21006799999999
$data[0,4] = 0.7601982364861632
$data[0,5] = 0.7601982364861634
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 117.044563
$data[0,9] = 351.133689
$data[0,10] = 0.3245365645427815
$data[0,11] = 0.3245365645427815
$data[0,12] = 1412.730528420095
$data[0,13] = 12714.57475578085
$data[0,14] = 0.2467121240407004
$data[0,15] = 0.2467121240407004
$data[1,0] = 3
$data[1,1] = 1
$data[1,2] = 12.07002266666666
$data[1,3] = 36.21006799999999
$data[1,4] = 0.7601982364861632
$data[1,5] = 0.7601982364861634
$data[1,6] = 3
$data[1,7] = 1
$data[1,8] = 101.5800373333333
$data[1,9] = 304.740112
$data[1,10] = 0.281657135515876
$data[1,11] = 0.281657135515876
$data[1,12] = 1226.073353094179
$data[1,13] = 11034.66017784761
$data[1,14] = 0.2141152577129132
$data[1,15] = 0.2141152577129132
$data[2,0] = 3
$data[2,1] = 1
$data[2,2] = 12.07002266666666
$data[2,3] = 36.21006799999999
$data[2,4] = 0.7601982364861632
$data[2,5] = 0.7601982364861634
$data[2,6] = 3
$data[2,7] = 1
$data[2,8] = 142.0267893333333
$data[2,9] = 426.080368
$data[2,10] = 0.3938062999413425
$data[2,11] = 0.3938062999413425
$data[2,12] = 1714.266566527225
$data[2,13] = 15428.39909874502
$data[2,14] = 0.2993708547325497
$data[2,15] = 0.2993708547325497
$data[3,0] = 3
$data[3,1] = 1
$data[3,2] = 1.308268
$data[3,3] = 3.924804
$data[3,4] = 0.08239777620284613
$data[3,5] = 0.08239777620284613
$data[3,6] = 3
$data[3,7] = 1
$data[3,8] = 117.044563
$data[3,9] = 351.133689
$data[3,10] = 0.3245365645427815
$data[3,11] = 0.3245365645427815
$data[3,12] = 153.125656346884
$data[3,13] = 1378.130907121956
$data[3,14] = 0.02674109121483664
$data[3,15] = 0.02674109121483664
$data[4,0] = 3
$data[4,1] = 1
$data[4,2] = 1.308268
$data[4,3] = 3.924804
$data[4,4] = 0.08239777620284613
$data[4,5] = 0.08239777620284613
$data[4,6] = 3
$data[4,7] = 1
$data[4,8] = 101.5800373333333
$data[4,9] = 304.740112
$data[4,10] = 0.281657135515876
$data[4,11] = 0.281657135515876
$data[4,12] = 132.8939122820053
$data[4,13] = 1196.045210538048
$data[4,14] = 0.02320792161817186
$data[4,15] = 0.02320792161817185
$data[5,0] = 3
$data[5,1] = 1
$data[5,2] = 1.308268
$data[5,3] = 3.924804
$data[5,4] = 0.08239777620284613
$data[5,5] = 0.08239777620284613
$data[5,6] = 3
$data[5,7] = 1
$data[5,8] = 142.0267893333333
$data[5,9] = 426.080368
$data[5,10] = 0.3938062999413425
$data[5,11] = 0.3938062999413425
$data[5,12] = 185.8091036275413
$data[5,13] = 1672.281932647872
$data[5,14] = 0.03244876336983764
$data[5,15] = 0.03244876336983764
$data[6,0] = 3
$data[6,1] = 1
$data[6,2] = 2.499176666666667
$data[6,3] = 7.49753
$data[6,4] = 0.1574039873109905
$data[6,5] = 0.1574039873109906
$data[6,6] = 3
$data[6,7] = 1
$data[6,8] = 117.044563
$data[6,9] = 351.133689
$data[6,10] = 0.3245365645427815
$data[6,11] = 0.3245365645427815
$data[6,12] = 292.5150408097966
$data[6,13] = 2632.63536728817
$data[6,14] = 0.05108334928724443
$data[6,15] = 0.05108334928724444
$data[7,0] = 3
$data[7,1] = 1
$data[7,2] = 2.499176666666667
$data[7,3] = 7.49753
$data[7,4] = 0.1574039873109905
$data[7,5] = 0.1574039873109906
$data[7,6] = 3
$data[7,7] = 1
$data[7,8] = 101.5800373333333
$data[7,9] = 304.740112
$data[7,10] = 0.281657135515876
$data[7,11] = 0.281657135515876
$data[7,12] = 253.8664591025956
$data[7,13] = 2284.79813192336
$data[7,14] = 0.04433395618479089
$data[7,15] = 0.04433395618479088
$data[8,0] = 3
$data[8,1] = 1
$data[8,2] = 2.499176666666667
$data[8,3] = 7.49753
$data[8,4] = 0.1574039873109905
$data[8,5] = 0.1574039873109906
$data[8,6] = 3
$data[8,7] = 1
$data[8,8] = 142.0267893333333
$data[8,9] = 426.080368
$data[8,10] = 0.3938062999413425
$data[8,11] = 0.3938062999413425
$data[8,12] = 354.9500379434489
$data[8,13] = 3194.55034149104
$data[8,14] = 0.06198668183895521
$data[8,15] = 0.06198668183895522

$ws.Range("E2:T10").Value = $data
